$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 120, pushing existing rows 120..175
# down to 122..177 (matches the target dimension A1:T177).
$ws.Rows("120:121").Insert()

# Populate the two newly inserted rows with the new weekly records.
# Row 120
$ws.Range("A120").Value = 3
$ws.Range("B120").Value = "Femacal de La Calera"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 44529
$ws.Range("E120").Value = 5
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100101
$ws.Range("H120").Value = "Berries"
$ws.Range("I120").Value = 100112025
$ws.Range("J120").Value = "Frutilla"
$ws.Range("K120").Value = "Sin especificar"
$ws.Range("L120").Value = "Especial"
$ws.Range("M120").Value = 85
$ws.Range("N120").Value = 6000
$ws.Range("O120").Value = 6000
$ws.Range("P120").Value = 6000
$ws.Range("Q120").Value = "$/bandeja 7 kilos"
$ws.Range("R120").Value = "Provincia de Melipilla"
$ws.Range("S120").Value = 857
$ws.Range("T120").Value = 7

# Row 121
$ws.Range("A121").Value = 3
$ws.Range("B121").Value = "Femacal de La Calera"
$ws.Range("C121").Value = "Coquimbo"
$ws.Range("D121").Value = 44529
$ws.Range("E121").Value = 5
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100101
$ws.Range("H121").Value = "Berries"
$ws.Range("I121").Value = 100112025
$ws.Range("J121").Value = "Frutilla"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Segunda"
$ws.Range("M121").Value = 70
$ws.Range("N121").Value = 4000
$ws.Range("O121").Value = 4000
$ws.Range("P121").Value = 4000
$ws.Range("Q121").Value = "$/bandeja 7 kilos"
$ws.Range("R121").Value = "Provincia de Melipilla"
$ws.Range("S121").Value = 571
$ws.Range("T121").Value = 7
